$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 0.04172184405617529
$ws.Range("C2").Value = 109.9114832445916
$ws.Range("D2").Value = 18.71679738969934
$bigval = 2.521694498980204 * [Math]::Pow(10, 27)
$ws.Range("E2").Value = $bigval
$ws.Range("G2").Value = $bigval
